$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '69.412.63'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -1.42%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.540.04'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -2.95%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '587.24'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +1.86%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '173.06'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -1.69%  '

# Row 7
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.62%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.536.70'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -2.84%  '

# Row 9
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.01%  '

# Row 10
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -3.53%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.79'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -0.36%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.584'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -3.34%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '47.56'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -2.12%  '

# Row 14
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -3.02%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.108.56'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -2.87%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '8.55'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -3.77%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '629.75'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -6.04%  '

# Row 18
$ws.Range('B18').NumberFormat = "@"
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').NumberFormat = "@"
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '69.447.58'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -1.50%  '

# Row 19
$ws.Range('B19').NumberFormat = "@"
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').NumberFormat = "@"
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.535.21'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -3.04%  '

# Row 20
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +1.38%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.46'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -1.71%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '11.20'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -1.56%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.892'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -4.15%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '15.99'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -6.61%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '97.52'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -2.85%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.83'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -1.71%  '

# Row 27
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +0.06%  '

# Row 28
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -4.94%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.35'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -6.32%  '

# Row 30
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -5.09%  '

# Row 31
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -4.13%  '

# Row 32
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -5.49%  '

# Row 33
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -3.68%  '

# Row 34
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -4.37%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '636.64'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +9.67%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '10.82'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -2.08%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.51'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -11.27%  '

# Row 38
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -3.39%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '57.35'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -1.22%  '

# Row 40
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +0.15%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0459'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +1.21%  '

# Row 42
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -2.93%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.395.83'

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.330'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -4.17%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0₃0705'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -5.18%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '32.92'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -5.52%  '

# Row 47
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -4.67%  '

# Row 48
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -6.56%  '

# Row 49
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -1.93%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '132.55'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -2.15%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '5.66'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +13.90%  '
